# Update various files to use previous month date
# This zeroes out the most-recent month's data column (oct_2020 = column H
# for rows 17,18,20,21,27,28,29 and nov_2020 = column I for rows
# 31,33,34,36,40,41,43,46,47) along with the corresponding SFY 2021 Total
# (column Q) for those same rows, matching the values previously reported
# for the month that has since rolled back to zero/unsubmitted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column H (oct_2020) and column Q (SFY 2021 Total) need to be zeroed
$hRows = @(17, 18, 20, 21, 27, 28, 29)
foreach ($r in $hRows) {
    $ws.Cells.Item($r, 8).Value = 0   # Column H
    $ws.Cells.Item($r, 17).Value = 0  # Column Q
}

# Rows where column I (nov_2020) and column Q (SFY 2021 Total) need to be zeroed
$iRows = @(31, 33, 34, 36, 40, 41, 43, 46, 47)
foreach ($r in $iRows) {
    $ws.Cells.Item($r, 9).Value = 0   # Column I
    $ws.Cells.Item($r, 17).Value = 0  # Column Q
}

# Selection moved to M18 as part of the resave
$ws.Range("M18").Select()
